$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.214.64"
$ws.Range("E2").Value = "  +2.61%  "

$ws.Range("D3").Value = "2.321.10"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.33%  "

$ws.Range("D9").Value = "2.318.18"
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").Value = "2.735.10"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "60.162.16"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").Value = "2.320.97"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("E25").Value = "  -0.84%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("E31").Value = "  +0.42%  "

$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("E34").Value = "  +8.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.380"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.24%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.563"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0492"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").Value = "0.0₆0229"
$ws.Range("E49").Value = "  +20.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.41%  "
